$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data added a new daily price observation for
# "Feria Lagunitas de Puerto Montt - Acelga" that belongs chronologically
# right after the existing row 186 record. In the canonical (row-ordered)
# layout this means: insert a new row at position 187 (pushing the old
# row 187 and everything below it down by one row), fill that new row
# with a duplicate of row 186's record, then give row 186 and the new
# row 187 their correct, distinct values.

# 1) Insert a blank row at 187; rows 187..282 shift down to 188..283.
$ws.Rows(187).Insert()

# 2) Seed the new row with a copy of row 186 (same market/category/etc.)
$ws.Range("A186:R186").Copy()
$ws.Range("A187").PasteSpecial()
$excel.CutCopyMode = $false

# 3) Row 186 becomes the new observation (new date + new volume).
$ws.Range("D186").Value = 44992
$ws.Range("J186").Value = 90

# Row 187 already carries row 186's former D/J values (44608 / 20) from
# the copy in step 2, so no further change is needed there.
